$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 727.8049
$ws.Range("J17").Value = 727.8049
$ws.Range("L17").Value = 2183.4147
$ws.Range("N17").Value = -2519.4147
$ws.Range("H19").Value = 1233.1666
$ws.Range("J19").Value = 1233.1666
$ws.Range("L19").Value = 1233.1666
$ws.Range("N19").Value = -1583.1666
$ws.Range("H54").Value = 21583.334
$ws.Range("I54").Value = 4750
$ws.Range("K54").Value = 4750
$ws.Range("M54").Value = -4264
$ws.Range("H125").Value = 888.35
$ws.Range("I125").Value = 214.5
$ws.Range("J125").Value = 1177.1428
$ws.Range("K125").Value = 1930.5
$ws.Range("L125").Value = 10594.2852
$ws.Range("M125").Value = 529.5
$ws.Range("N125").Value = -15514.2852
$ws.Range("H131").Value = 996.1111
$ws.Range("I131").Value = 840.7143
$ws.Range("J131").Value = 1050.5
$ws.Range("K131").Value = 2522.1429
$ws.Range("L131").Value = 3151.5
$ws.Range("M131").Value = 2517.8571
$ws.Range("N131").Value = -13231.5
$ws.Range("H132").Value = 792351.7
$ws.Range("I132").Value = 1929.5769
$ws.Range("J132").Value = 4902546.5
$ws.Range("K132").Value = 5788.7307
$ws.Range("L132").Value = 14707639.5
$ws.Range("M132").Value = -3258.7307
$ws.Range("N132").Value = -14712699.5
$ws.Range("H137").Value = 1668036.1
$ws.Range("I137").Value = 2326588.8
$ws.Range("J137").Value = 2285.353
$ws.Range("K137").Value = 6979766.399999999
$ws.Range("L137").Value = 6856.059
$ws.Range("M137").Value = -6977216.399999999
$ws.Range("N137").Value = -11956.059
$ws.Range("H138").Value = 2223686.8
$ws.Range("I138").Value = 1089.3492
$ws.Range("J138").Value = 13892323
$ws.Range("K138").Value = 3268.0476
$ws.Range("L138").Value = 41676969
$ws.Range("M138").Value = 1871.9524
$ws.Range("N138").Value = -41687249
$ws.Range("H141").Value = 1817.3704
$ws.Range("I141").Value = 1234.7446
$ws.Range("K141").Value = 3704.2338
$ws.Range("M141").Value = 1475.7662
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 17545526
$ws.Range("I5").Value = 26317288
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 26317288
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = -26317176
$ws.Range("N5").Value = -2224
$ws.Range("H122").Value = 3474126.5
$ws.Range("I122").Value = 1960.2593
$ws.Range("K122").Value = 5880.7779
$ws.Range("M122").Value = -3430.7779
$ws.Range("H132").Value = 43351.855
$ws.Range("I132").Value = 25643.35
$ws.Range("K132").Value = 76930.04999999999
$ws.Range("M132").Value = -74400.04999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 17545526
$ws.Range("I4").Value = 26317288
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 26317288
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = -26317173
$ws.Range("N4").Value = -2230
$ws.Range("H105").Value = 20835710
$ws.Range("I105").Value = 38463870
$ws.Range("K105").Value = 38463870
$ws.Range("M105").Value = -38462123
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3258.0588
$ws.Range("I31").Value = 1272.3077
$ws.Range("J31").Value = 9711.75
$ws.Range("K31").Value = 1272.3077
$ws.Range("L31").Value = 9711.75
$ws.Range("M31").Value = -977.3077000000001
$ws.Range("N31").Value = -10301.75
$ws.Range("H34").Value = 3258.0588
$ws.Range("I34").Value = 1272.3077
$ws.Range("J34").Value = 9711.75
$ws.Range("K34").Value = 1272.3077
$ws.Range("L34").Value = 9711.75
$ws.Range("M34").Value = -1070.3077
$ws.Range("N34").Value = -10115.75
$ws.Range("H58").Value = 17858542
$ws.Range("I58").Value = 23257208
$ws.Range("J58").Value = 1415
$ws.Range("K58").Value = 23257208
$ws.Range("L58").Value = 1415
$ws.Range("M58").Value = -23257005
$ws.Range("N58").Value = -1821
$ws.Range("H132").Value = 16815.562
$ws.Range("I132").Value = 1111.4423
$ws.Range("J132").Value = 84866.75
$ws.Range("K132").Value = 3334.3269
$ws.Range("L132").Value = 254600.25
$ws.Range("M132").Value = -804.3269
$ws.Range("N132").Value = -259660.25
$ws.Range("H136").Value = 17858542
$ws.Range("I136").Value = 23257208
$ws.Range("J136").Value = 1415
$ws.Range("K136").Value = 69771624
$ws.Range("L136").Value = 4245
$ws.Range("M136").Value = -69769074
$ws.Range("N136").Value = -9345
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 559.9231
$ws.Range("I5").Value = 458.5
$ws.Range("J5").Value = 788.125
$ws.Range("K5").Value = 1375.5
$ws.Range("L5").Value = 2364.375
$ws.Range("M5").Value = -1263.5
$ws.Range("N5").Value = -2588.375
$ws.Range("H122").Value = 509.6389
$ws.Range("I122").Value = 256.68
$ws.Range("J122").Value = 1084.5454
$ws.Range("K122").Value = 2310.12
$ws.Range("L122").Value = 9760.908599999999
$ws.Range("M122").Value = 139.8800000000001
$ws.Range("N122").Value = -14660.9086
$ws.Range("H131").Value = 1120.7115
$ws.Range("I131").Value = 627.8
$ws.Range("J131").Value = 1173.1489
$ws.Range("K131").Value = 1883.4
$ws.Range("L131").Value = 3519.4467
$ws.Range("M131").Value = 3156.6
$ws.Range("N131").Value = -13599.4467
$ws.Range("H132").Value = 2844.7222
$ws.Range("I132").Value = 1762.5
$ws.Range("J132").Value = 3710.5
$ws.Range("K132").Value = 15862.5
$ws.Range("L132").Value = 33394.5
$ws.Range("M132").Value = -13332.5
$ws.Range("N132").Value = -38454.5
$ws.Range("H135").Value = 559.9231
$ws.Range("I135").Value = 458.5
$ws.Range("J135").Value = 788.125
$ws.Range("K135").Value = 4126.5
$ws.Range("L135").Value = 7093.125
$ws.Range("M135").Value = -1591.5
$ws.Range("N135").Value = -12163.125
$ws.Range("H139").Value = 4467.6875
$ws.Range("I139").Value = 2422.6365
$ws.Range("J139").Value = 5075.676
$ws.Range("K139").Value = 7267.9095
$ws.Range("L139").Value = 15227.028
$ws.Range("M139").Value = -2127.9095
$ws.Range("N139").Value = -25507.028
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 850
$ws.Range("J122").Value = 1933.3334
$ws.Range("K122").Value = 2550
$ws.Range("L122").Value = 5800.0002
$ws.Range("M122").Value = -100
$ws.Range("N122").Value = -10700.0002
$ws.Range("H123").Value = 29123.363
$ws.Range("J123").Value = 29123.363
$ws.Range("L123").Value = 29123.363
$ws.Range("N123").Value = -34023.363
$ws.Range("H132").Value = 43686.297
$ws.Range("I132").Value = 26665.154
$ws.Range("J132").Value = 126664.375
$ws.Range("K132").Value = 79995.462
$ws.Range("L132").Value = 379993.125
$ws.Range("M132").Value = -77465.462
$ws.Range("N132").Value = -385053.125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2954.3684
$ws.Range("I122").Value = 2547.65
$ws.Range("J122").Value = 3406.2778
$ws.Range("K122").Value = 7642.950000000001
$ws.Range("L122").Value = 10218.8334
$ws.Range("M122").Value = -5192.950000000001
$ws.Range("N122").Value = -15118.8334
$ws.Range("H136").Value = 82705.39999999999
$ws.Range("I136").Value = 50776.75
$ws.Range("K136").Value = 152330.25
$ws.Range("M136").Value = -149780.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1684.25
$ws.Range("I122").Value = 1210.1154
$ws.Range("K122").Value = 3630.3462
$ws.Range("M122").Value = -1180.3462
$ws.Range("H123").Value = 33517.617
$ws.Range("J123").Value = 33517.617
$ws.Range("L123").Value = 33517.617
$ws.Range("N123").Value = -43317.617
$ws.Range("H132").Value = 39165.145
$ws.Range("I132").Value = 23411.797
$ws.Range("K132").Value = 70235.391
$ws.Range("M132").Value = -67705.391
$ws.Range("H136").Value = 47455.977
$ws.Range("I136").Value = 36538.855
$ws.Range("J136").Value = 66560.94
$ws.Range("K136").Value = 109616.565
$ws.Range("L136").Value = 199682.82
$ws.Range("M136").Value = -107066.565
$ws.Range("N136").Value = -204782.82
